$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 4 into a new row 5 (same formatting, same values), then
# change the "Opprotunity Name" (column D) value to the new test name.
$ws.Range("A4:P4").Copy($ws.Range("A5:P5"))
$ws.Range("D5").Value = "Test Second Automation"

# Widen column D to fit the longer text.
$ws.Columns("D").ColumnWidth = 20.43

# Give A5 the same mailto hyperlink as A4.
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:APPTESTMBOB4@netapp.com")

# Adding the hyperlink reset A5's cell style (Excel applies a fresh
# Hyperlink-style xf); restore the original Hyperlink-style formatting
# used by A4 so A5 matches it exactly.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
